# The uploaded workbook now has its data shifted one column to the left:
# what used to live in B1:E8 (headers "Longitud2".."Longitud5" plus the
# measurements) now lives in A1:D8. Deleting the (empty) column A achieves
# exactly that shift while letting Excel naturally re-flow the column
# styles/widths that used to be attached to B:E onto A:D.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("A").Delete()

# Leave the selection at A1, matching the default (no stale B1:E1048576
# selection left over from when the data used to start in column B).
$ws.Range("A1").Select() | Out-Null
